# Update the date heading at the top of the document.
$d = $word.ActiveDocument
[void]$d.Content.Find.Execute("2026-01-08 Thursday", $true, $false, $false, $false, $false,
                               $true, 1, $false, "2026-01-09 Friday", 2)

# Update the division problems in the table. Addressed by (row, column)
# because several cells share identical text ("59÷2=29, 1" appears twice)
# and map to different replacements, so a blind Find/Replace would be unsafe.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "40÷6=6, 4"
$t.Cell(1, 2).Range.Text = "42÷6=7, 0"
$t.Cell(1, 3).Range.Text = "66÷9=7, 3"
$t.Cell(1, 4).Range.Text = "40÷5=8, 0"
$t.Cell(1, 5).Range.Text = "63÷2=31, 1"

$t.Cell(5, 1).Range.Text = "25÷2=12, 1"
$t.Cell(5, 2).Range.Text = "73÷7=10, 3"
$t.Cell(5, 3).Range.Text = "32÷2=16, 0"
$t.Cell(5, 4).Range.Text = "60÷3=20, 0"
$t.Cell(5, 5).Range.Text = "60÷9=6, 6"

$t.Cell(9, 1).Range.Text = "55÷4=13, 3"
$t.Cell(9, 2).Range.Text = "16÷7=2, 2"
$t.Cell(9, 3).Range.Text = "39÷3=13, 0"
$t.Cell(9, 4).Range.Text = "25÷5=5, 0"
$t.Cell(9, 5).Range.Text = "69÷9=7, 6"

$t.Cell(13, 1).Range.Text = "64÷8=8, 0"
$t.Cell(13, 2).Range.Text = "42÷6=7, 0"
$t.Cell(13, 3).Range.Text = "79÷4=19, 3"
$t.Cell(13, 4).Range.Text = "86÷6=14, 2"
$t.Cell(13, 5).Range.Text = "23÷6=3, 5"

$t.Cell(17, 1).Range.Text = "70÷7=10, 0"
$t.Cell(17, 2).Range.Text = "33÷7=4, 5"
$t.Cell(17, 3).Range.Text = "85÷8=10, 5"
$t.Cell(17, 4).Range.Text = "91÷7=13, 0"
$t.Cell(17, 5).Range.Text = "29÷6=4, 5"
